# Dev IV Project Rubric.xlsx - edit script
# Commit: "Got all objects drawn and setup multiple viewports. Have point
# light code setup in pixel shader and main. Further work needed to get
# points."
#
# The underlying change is the grader marking several milestone rows as
# "III" (Student column) with "X" (Confidence column) on Sheet1, plus two
# of the carry-over rows (86/87) getting an "X" in the Student column.
# There is also a cosmetic change of the saved sheet view (scroll
# position / zoom / selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark these rows as Milestone "III" achieved, with confidence confirmed ("X").
# (G column formulas auto-recalculate from these inputs.)
$rows = @(5, 7, 8, 9, 68, 79)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "III"
    $ws.Range("F$r").Value = "X"
}

# Carry-over rows: mark Student column with "X" as well.
$ws.Range("E86").Value = "X"
$ws.Range("E87").Value = "X"

# View state: scrolled/zoomed/selected differently after the edit session.
$excel.ActiveWindow.Zoom = 90
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F18").Select()
